# Auto-generated edit script: updates cached price/profit values per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 10253.6
$ws.Range("I9").Value = 12704.5
$ws.Range("K9").Value = 12704.5
$ws.Range("M9").Value = -12535.5
$ws.Range("H19").Value = 1455.2
$ws.Range("I19").Value = 1500
$ws.Range("K19").Value = 1500
$ws.Range("M19").Value = -1325
$ws.Range("H41").Value = 931.4666999999999
$ws.Range("I41").Value = 1353
$ws.Range("J41").Value = 88.40000000000001
$ws.Range("K41").Value = 1353
$ws.Range("L41").Value = 88.40000000000001
$ws.Range("M41").Value = -913
$ws.Range("N41").Value = -968.4
$ws.Range("H88").Value = 9148.777
$ws.Range("I88").Value = 8156
$ws.Range("J88").Value = 9347.333000000001
$ws.Range("K88").Value = 8156
$ws.Range("L88").Value = 9347.333000000001
$ws.Range("M88").Value = -7750
$ws.Range("N88").Value = -10159.333
$ws.Range("H91").Value = 9148.777
$ws.Range("I91").Value = 8156
$ws.Range("J91").Value = 9347.333000000001
$ws.Range("K91").Value = 8156
$ws.Range("L91").Value = 9347.333000000001
$ws.Range("M91").Value = -6752
$ws.Range("N91").Value = -12155.333
$ws.Range("H103").Value = 808.8333
$ws.Range("J103").Value = 573.7778
$ws.Range("L103").Value = 1721.3334
$ws.Range("N103").Value = -2893.3334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1052429.4
$ws.Range("I2").Value = 2103103.2
$ws.Range("J2").Value = 1755.5714
$ws.Range("K2").Value = 2103103.2
$ws.Range("L2").Value = 1755.5714
$ws.Range("M2").Value = -2102990.2
$ws.Range("N2").Value = -1981.5714
$ws.Range("H32").Value = 2888.453
$ws.Range("I32").Value = 1185.2041
$ws.Range("K32").Value = 1185.2041
$ws.Range("M32").Value = -898.2040999999999
$ws.Range("H45").Value = 4335.304
$ws.Range("I45").Value = 5047
$ws.Range("J45").Value = 2318.8333
$ws.Range("K45").Value = 5047
$ws.Range("L45").Value = 2318.8333
$ws.Range("M45").Value = -4670
$ws.Range("N45").Value = -3072.8333
$ws.Range("H97").Value = 846.8333
$ws.Range("I97").Value = 832.9091
$ws.Range("K97").Value = 832.9091
$ws.Range("M97").Value = -336.9091
$ws.Range("H116").Value = 1052429.4
$ws.Range("I116").Value = 2103103.2
$ws.Range("J116").Value = 1755.5714
$ws.Range("K116").Value = 2103103.2
$ws.Range("L116").Value = 1755.5714
$ws.Range("M116").Value = -2100809.2
$ws.Range("N116").Value = -6343.5714
$ws.Range("H132").Value = 3228359.2
$ws.Range("I132").Value = 3228359.2
$ws.Range("K132").Value = 9685077.600000001
$ws.Range("M132").Value = -9682547.600000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1052429.4
$ws.Range("I3").Value = 2103103.2
$ws.Range("J3").Value = 1755.5714
$ws.Range("K3").Value = 2103103.2
$ws.Range("L3").Value = 1755.5714
$ws.Range("M3").Value = -2102989.2
$ws.Range("N3").Value = -1983.5714
$ws.Range("H47").Value = 495000
$ws.Range("J47").Value = 495000
$ws.Range("L47").Value = 495000
$ws.Range("N47").Value = -496040
$ws.Range("H75").Value = 15490.333
$ws.Range("J75").Value = 17235.5
$ws.Range("L75").Value = 17235.5
$ws.Range("N75").Value = -19107.5
$ws.Range("H78").Value = 15490.333
$ws.Range("J78").Value = 17235.5
$ws.Range("L78").Value = 51706.5
$ws.Range("N78").Value = -61066.5
$ws.Range("H94").Value = 2358.75
$ws.Range("I94").Value = 2160
$ws.Range("K94").Value = 2160
$ws.Range("M94").Value = -1709
$ws.Range("H99").Value = 2007.7646
$ws.Range("I99").Value = 1985.8889
$ws.Range("J99").Value = 2032.375
$ws.Range("K99").Value = 1985.8889
$ws.Range("L99").Value = 2032.375
$ws.Range("M99").Value = -487.8888999999999
$ws.Range("N99").Value = -5028.375
$ws.Range("H100").Value = 22939.834
$ws.Range("J100").Value = 22939.834
$ws.Range("L100").Value = 22939.834
$ws.Range("N100").Value = -25103.834
$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988
$ws.Range("H105").Value = 2522.2856
$ws.Range("I105").Value = 2042.0834
$ws.Range("K105").Value = 2042.0834
$ws.Range("M105").Value = -295.0834

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2999.5
$ws.Range("I62").Value = 2999.5
$ws.Range("K62").Value = 2999.5
$ws.Range("M62").Value = -2375.5
$ws.Range("H65").Value = 2999.5
$ws.Range("I65").Value = 2999.5
$ws.Range("K65").Value = 14997.5
$ws.Range("M65").Value = -11877.5
$ws.Range("H105").Value = 1857388.5
$ws.Range("I105").Value = 2269519.2
$ws.Range("K105").Value = 2269519.2
$ws.Range("M105").Value = -2267772.2
$ws.Range("H120").Value = 59998.75
$ws.Range("J120").Value = 59998.75
$ws.Range("L120").Value = 59998.75
$ws.Range("N120").Value = -67256.75
$ws.Range("H122").Value = 1776.4445
$ws.Range("J122").Value = 1094
$ws.Range("L122").Value = 3282
$ws.Range("N122").Value = -8182
$ws.Range("H134").Value = 7385978
$ws.Range("I134").Value = 8370227.5
$ws.Range("K134").Value = 25110682.5
$ws.Range("M134").Value = -25108147.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1432233
$ws.Range("I7").Value = 5004997.5
$ws.Range("J7").Value = 3127.2
$ws.Range("K7").Value = 15014992.5
$ws.Range("L7").Value = 9381.599999999999
$ws.Range("M7").Value = -15014880.5
$ws.Range("N7").Value = -9605.599999999999
$ws.Range("H23").Value = 413.8
$ws.Range("J23").Value = 489.75
$ws.Range("L23").Value = 1469.25
$ws.Range("N23").Value = -1939.25
$ws.Range("H57").Value = 20000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 20000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 60000
$ws.Range("M57").Value = $null
$ws.Range("N57").Value = -61118
$ws.Range("H117").Value = 1874
$ws.Range("I117").Value = 355.6
$ws.Range("J117").Value = 2717.5557
$ws.Range("K117").Value = 1066.8
$ws.Range("L117").Value = 8152.6671
$ws.Range("M117").Value = 2375.2
$ws.Range("N117").Value = -15036.6671
$ws.Range("H125").Value = 8000
$ws.Range("I125").Value = 8000
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 24000
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -19080
$ws.Range("N125").Value = $null
$ws.Range("H140").Value = 975.44446
$ws.Range("I140").Value = 975.44446
$ws.Range("K140").Value = 2926.33338
$ws.Range("M140").Value = 2253.66662

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.57143000000001
$ws.Range("I2").Value = 34.375
$ws.Range("K2").Value = 34.375
$ws.Range("M2").Value = 78.625
$ws.Range("H122").Value = 68397.72
$ws.Range("I122").Value = 68397.72
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 205193.16
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -202743.16
$ws.Range("N122").Value = $null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 39999
$ws.Range("J99").Value = 39999
$ws.Range("L99").Value = 39999
$ws.Range("N99").Value = -45989
$ws.Range("H102").Value = 75280
$ws.Range("J102").Value = 75280
$ws.Range("L102").Value = 75280
$ws.Range("N102").Value = -81770
$ws.Range("H132").Value = 9620781
$ws.Range("I132").Value = 10422347
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 31267041
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -31264511
$ws.Range("N132").Value = -11058.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 97999.60000000001
$ws.Range("I75").Value = 35000
$ws.Range("K75").Value = 35000
$ws.Range("M75").Value = -34064
$ws.Range("H78").Value = 97999.60000000001
$ws.Range("I78").Value = 35000
$ws.Range("K78").Value = 105000
$ws.Range("M78").Value = -100320
$ws.Range("H102").Value = 68499.5
$ws.Range("J102").Value = 68499.5
$ws.Range("L102").Value = 68499.5
$ws.Range("N102").Value = -74989.5
$ws.Range("H130").Value = 49999
$ws.Range("J130").Value = 49999
$ws.Range("L130").Value = 49999
$ws.Range("N130").Value = -60039
$ws.Range("H132").Value = 13162702
$ws.Range("I132").Value = 17859156
$ws.Range("K132").Value = 53577468
$ws.Range("M132").Value = -53574938
